# Comercializadora del Agro de Limarí - Chirimoya
# Add a new weekly price report (bandeja de 10 kilos) as the most recent
# entries, inserted above the existing historical rows (which simply get
# pushed down by 3 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at the top of the data block (before the old
# row 44), pushing all existing data rows down by 3.
$ws.Rows("44:46").Insert()

# New row 44: Especial
$ws.Cells.Item(44, 1).Value = 2
$ws.Cells.Item(44, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(44, 3).Value = "Coquimbo"
$ws.Cells.Item(44, 4).Value = 44818
$ws.Cells.Item(44, 5).Value = 4
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100107
$ws.Cells.Item(44, 8).Value = "Otros"
$ws.Cells.Item(44, 9).Value = 100107002
$ws.Cells.Item(44, 10).Value = "Chirimoya"
$ws.Cells.Item(44, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(44, 12).Value = "Especial"
$ws.Cells.Item(44, 13).Value = 240
$ws.Cells.Item(44, 14).Value = 23000
$ws.Cells.Item(44, 15).Value = 24000
$ws.Cells.Item(44, 16).Value = 23500
$ws.Cells.Item(44, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(44, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(44, 19).Value = 2350
$ws.Cells.Item(44, 20).Value = 10

# New row 45: Primera
$ws.Cells.Item(45, 1).Value = 2
$ws.Cells.Item(45, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(45, 3).Value = "Coquimbo"
$ws.Cells.Item(45, 4).Value = 44818
$ws.Cells.Item(45, 5).Value = 4
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100107
$ws.Cells.Item(45, 8).Value = "Otros"
$ws.Cells.Item(45, 9).Value = 100107002
$ws.Cells.Item(45, 10).Value = "Chirimoya"
$ws.Cells.Item(45, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 300
$ws.Cells.Item(45, 14).Value = 19000
$ws.Cells.Item(45, 15).Value = 20000
$ws.Cells.Item(45, 16).Value = 19500
$ws.Cells.Item(45, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(45, 19).Value = 1950
$ws.Cells.Item(45, 20).Value = 10

# New row 46: Segunda
$ws.Cells.Item(46, 1).Value = 2
$ws.Cells.Item(46, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(46, 3).Value = "Coquimbo"
$ws.Cells.Item(46, 4).Value = 44818
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100107
$ws.Cells.Item(46, 8).Value = "Otros"
$ws.Cells.Item(46, 9).Value = 100107002
$ws.Cells.Item(46, 10).Value = "Chirimoya"
$ws.Cells.Item(46, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(46, 12).Value = "Segunda"
$ws.Cells.Item(46, 13).Value = 240
$ws.Cells.Item(46, 14).Value = 16000
$ws.Cells.Item(46, 15).Value = 17000
$ws.Cells.Item(46, 16).Value = 16500
$ws.Cells.Item(46, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 19).Value = 1650
$ws.Cells.Item(46, 20).Value = 10
